$d = $word.ActiveDocument

# 1. Fix date typo: "Sep 2023 - Present" -> "Sep 2020 - Jun 2023" (Betterhealth job)
$d.Content.Find.Execute("Sep 2023 " + [char]0x2013 + " Present", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sep 2020 " + [char]0x2013 + " Jun 2023", 2)

Write-Output "done"
